# Update "想去人数" (F column) values across the four sheets of the workbook.
$wb = $excel.ActiveWorkbook

# 展览 (sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 2498
$ws1.Range("F13").Value = 1472
$ws1.Range("F14").Value = 1472
$ws1.Range("F15").Value = 1221
$ws1.Range("F17").Value = 3555
$ws1.Range("F19").Value = 3276
$ws1.Range("F20").Value = 734
$ws1.Range("F21").Value = 2122
$ws1.Range("F29").Value = 946

# 演出 (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 206

# 本地生活 (sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 496

# 全部类型 (sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F11").Value = 496
$ws4.Range("F12").Value = 2498
$ws4.Range("F25").Value = 1472
$ws4.Range("F26").Value = 1472
$ws4.Range("F28").Value = 206
$ws4.Range("F29").Value = 1221
$ws4.Range("F32").Value = 3555
$ws4.Range("F34").Value = 3276
$ws4.Range("F35").Value = 734
$ws4.Range("F37").Value = 2122
$ws4.Range("F49").Value = 946
